$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4101
$ws.Range("J88").Value = 4451.5
$ws.Range("L88").Value = 4451.5
$ws.Range("N88").Value = -5263.5
$ws.Range("H91").Value = 4101
$ws.Range("J91").Value = 4451.5
$ws.Range("L91").Value = 4451.5
$ws.Range("N91").Value = -7259.5
$ws.Range("H98").Value = 3463.6
$ws.Range("I98").Value = 3463.6
$ws.Range("K98").Value = 3463.6
$ws.Range("M98").Value = -1965.6
$ws.Range("H122").Value = 3463.6
$ws.Range("I122").Value = 3463.6
$ws.Range("K122").Value = 10390.8
$ws.Range("M122").Value = -7940.799999999999
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 7699.8
$ws.Range("I127").Value = 7374.75
$ws.Range("J127").Value = 9000
$ws.Range("K127").Value = 22124.25
$ws.Range("L127").Value = 27000
$ws.Range("M127").Value = -17164.25
$ws.Range("N127").Value = -36920
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 3466.6667
$ws.Range("I129").Value = 1200
$ws.Range("J129").Value = 8000
$ws.Range("K129").Value = 3600
$ws.Range("L129").Value = 24000
$ws.Range("M129").Value = 1400
$ws.Range("N129").Value = -34000
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 250
$ws.Range("I131").Value = 250
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 750
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 4290
$ws.Range("H132").Value = 2708.0715
$ws.Range("I132").Value = 2708.0715
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8124.2145
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5594.2145
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 70707
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 70707
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 70707
$ws.Range("N134").Value = -80847
$ws.Range("H135").Value = 1170.3077
$ws.Range("I135").Value = 929.36365
$ws.Range("J135").Value = 2495.5
$ws.Range("K135").Value = 8364.272849999999
$ws.Range("L135").Value = 22459.5
$ws.Range("M135").Value = -5829.272849999999
$ws.Range("N135").Value = -27529.5
$ws.Range("H136").Value = 90000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 90000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -100200
$ws.Range("H137").Value = 11115
$ws.Range("I137").Value = 10384
$ws.Range("J137").Value = 12333.333
$ws.Range("K137").Value = 31152
$ws.Range("L137").Value = 36999.999
$ws.Range("M137").Value = -28602
$ws.Range("N137").Value = -42099.999
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 781.4286
$ws.Range("I141").Value = 781.4286
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2344.2858
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2835.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 4999
$ws.Range("I17").Value = 4999
$ws.Range("K17").Value = 4999
$ws.Range("M17").Value = -4826
$ws.Range("H23").Value = 5006
$ws.Range("J23").Value = 5006
$ws.Range("L23").Value = 5006
$ws.Range("N23").Value = -5524
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 45000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -55200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 571.25
$ws.Range("I107").Value = 526.8333
$ws.Range("K107").Value = 526.8333
$ws.Range("M107").Value = 1393.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 420.6
$ws.Range("J22").Value = 750
$ws.Range("L22").Value = 750
$ws.Range("N22").Value = -1450
$ws.Range("H115").Value = 59775
$ws.Range("J115").Value = 59775
$ws.Range("L115").Value = 59775
$ws.Range("N115").Value = -62125
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 46599.668
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 46599.668
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 46599.668
$ws.Range("N131").Value = -56679.668
$ws.Range("H132").Value = 9197.6
$ws.Range("I132").Value = 6662.6665
$ws.Range("J132").Value = 13000
$ws.Range("K132").Value = 19987.9995
$ws.Range("L132").Value = 39000
$ws.Range("M132").Value = -17457.9995
$ws.Range("N132").Value = -44060
$ws.Range("H133").Value = 50000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H134").Value = 4932.7856
$ws.Range("I134").Value = 2328.3333
$ws.Range("J134").Value = 9620.799999999999
$ws.Range("K134").Value = 6984.999899999999
$ws.Range("L134").Value = 28862.4
$ws.Range("M134").Value = -4449.999899999999
$ws.Range("N134").Value = -33932.39999999999
$ws.Range("H135").Value = 50000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 50000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 104824
$ws.Range("I141").Value = 90296
$ws.Range("J141").Value = 109666.664
$ws.Range("K141").Value = 90296
$ws.Range("L141").Value = 109666.664
$ws.Range("M141").Value = -85116
$ws.Range("N141").Value = -120026.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 848.1667
$ws.Range("I97").Value = 733
$ws.Range("J97").Value = 963.3333
$ws.Range("K97").Value = 733
$ws.Range("L97").Value = 963.3333
$ws.Range("M97").Value = -237
$ws.Range("N97").Value = -1955.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1950
$ws.Range("I82").Value = 1950
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1950
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -1589
$ws.Range("H85").Value = 1950
$ws.Range("I85").Value = 1950
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1950
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -702

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4662.6665
$ws.Range("I2").Value = 4662.6665
$ws.Range("K2").Value = 4662.6665
$ws.Range("M2").Value = -4550.6665

Write-Output "edit complete"